$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "51.041.97"
$ws.Cells.Item(2, 5).Value = "  -1.66%  "
$ws.Cells.Item(3, 4).Value = "2.943.96"
$ws.Cells.Item(3, 5).Value = "  -2.19%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "375.99"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -2.51%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "101.24"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -3.88%  "
$ws.Cells.Item(7, 5).Value = "  -1.78%  "
$ws.Cells.Item(8, 5).Value = "  +0.00%  "
$ws.Cells.Item(9, 5).Value = "  -1.82%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "36.42"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -3.40%  "
$ws.Cells.Item(11, 5).Value = "  -0.76%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.0852"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.03%  "
$ws.Cells.Item(13, 4).Value = "3.403.48"
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "18.16"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -2.20%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "7.62"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.66%  "
$ws.Cells.Item(16, 4).Value = "2.939.57"
$ws.Cells.Item(16, 5).Value = "  -1.97%  "
$ws.Cells.Item(17, 5).Value = "  -2.98%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "11.10"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +48.03%  "
$ws.Cells.Item(19, 4).Value = "51.019.57"
$ws.Cells.Item(19, 5).Value = "  -1.57%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "3.09"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -7.28%  "
$ws.Cells.Item(21, 5).Value = "  -4.60%  "
$ws.Cells.Item(22, 5).Value = "  -1.33%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "265.30"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.26%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "68.81"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.75%  "
$ws.Cells.Item(25, 5).Value = "  +6.76%  "
$ws.Cells.Item(26, 5).Value = "  -3.08%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "7.55"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -2.04%  "
$ws.Cells.Item(28, 5).Value = "  +0.05%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "25.69"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -2.09%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "0.164"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -4.86%  "
$ws.Cells.Item(31, 5).Value = "  -5.22%  "
$ws.Cells.Item(32, 5).Value = "  +0.50%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "50.86"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.74%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "2.05"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.94%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "33.51"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -5.12%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.0443"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -3.16%  "
$ws.Cells.Item(37, 5).Value = "  -0.06%  "
$ws.Cells.Item(38, 5).Value = "  +2.88%  "
$ws.Cells.Item(39, 5).Value = "  -1.16%  "
$ws.Cells.Item(40, 5).Value = "  -4.28%  "
$ws.Cells.Item(41, 5).Value = "  -3.42%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "2.49"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -5.55%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "120.38"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.57%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "21.35"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -3.12%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "3.39"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.65%  "
$ws.Cells.Item(46, 2).Value = "WEMIXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "2.03"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.31%  "
$ws.Cells.Item(47, 2).Value = "TheGraph"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.273"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -3.02%  "
$ws.Cells.Item(48, 5).Value = "  -2.12%  "
$ws.Cells.Item(49, 4).Value = "1.991.27"
$ws.Cells.Item(49, 5).Value = "  -2.83%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.0330"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -2.12%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "1.32"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.64%  "
